$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text changed from "Ready for handoff" to "Handback transform failed"
# for the row describing 5d98d244-57b3-486e-9c83-62eb5e34c1df across the
# Overview, zh-cn and de-de sheets.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# New "Error Detail" text added for the handback mismatch on both language sheets.
$wsZhCn.Range("K3").Value = "Handback file name: akd14ahf.aca is different with handoff file name: 5d98d244-57b3-486e-9c83-62eb5e34c1df.6709ea68443c7975dabe17629c1712f673b17c39.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: akd14ahf.aca is different with handoff file name: 5d98d244-57b3-486e-9c83-62eb5e34c1df.6709ea68443c7975dabe17629c1712f673b17c39.de-de."
